$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.603177
$ws.Range("H2").Value = 19.809531
$ws.Range("I2").Value = 0.5135477412645301
$ws.Range("J2").Value = 0.5135477412645302
$ws.Range("M2").Value = 43.73434833333334
$ws.Range("N2").Value = 131.203045
$ws.Range("O2").Value = 0.1998633389969613
$ws.Range("P2").Value = 0.1998633389969613
$ws.Range("Q2").Value = 288.785643024655
$ws.Range("R2").Value = 2599.070787221895
$ws.Range("S2").Value = 0.1026393663034766
$ws.Range("T2").Value = 0.1026393663034765
$ws.Range("G3").Value = 6.603177
$ws.Range("H3").Value = 19.809531
$ws.Range("I3").Value = 0.5135477412645301
$ws.Range("J3").Value = 0.5135477412645302
$ws.Range("O3").Value = 0.3183113588032023
$ws.Range("P3").Value = 0.3183113588032022
$ws.Range("Q3").Value = 459.933026713979
$ws.Range("R3").Value = 4139.397240425811
$ws.Range("S3").Value = 0.1634680793322279
$ws.Range("T3").Value = 0.1634680793322279
$ws.Range("G4").Value = 6.603177
$ws.Range("H4").Value = 19.809531
$ws.Range("I4").Value = 0.5135477412645301
$ws.Range("J4").Value = 0.5135477412645302
$ws.Range("M4").Value = 37.39234266666667
$ws.Range("N4").Value = 112.177028
$ws.Range("O4").Value = 0.1708807549004341
$ws.Range("P4").Value = 0.170880754900434
$ws.Range("Q4").Value = 246.908257072652
$ws.Range("R4").Value = 2222.174313653868
$ws.Range("S4").Value = 0.08775542570469569
$ws.Range("T4").Value = 0.08775542570469569
$ws.Range("G5").Value = 6.603177
$ws.Range("H5").Value = 19.809531
$ws.Range("I5").Value = 0.5135477412645301
$ws.Range("J5").Value = 0.5135477412645302
$ws.Range("M5").Value = 68.04127866666666
$ws.Range("N5").Value = 204.123836
$ws.Range("O5").Value = 0.3109445472994024
$ws.Range("P5").Value = 0.3109445472994024
$ws.Range("Q5").Value = 449.2886063423239
$ws.Range("R5").Value = 4043.597457080916
$ws.Range("S5").Value = 0.1596848699241299
$ws.Range("T5").Value = 0.15968486992413
$ws.Range("I6").Value = 0.02944398858046029
$ws.Range("J6").Value = 0.0294439885804603
$ws.Range("M6").Value = 43.73434833333334
$ws.Range("N6").Value = 131.203045
$ws.Range("O6").Value = 0.1998633389969613
$ws.Range("P6").Value = 0.1998633389969613
$ws.Range("Q6").Value = 16.55737235740056
$ws.Range("R6").Value = 149.016351216605
$ws.Range("S6").Value = 0.005884773871079193
$ws.Range("T6").Value = 0.005884773871079193
$ws.Range("I7").Value = 0.02944398858046029
$ws.Range("J7").Value = 0.0294439885804603
$ws.Range("O7").Value = 0.3183113588032023
$ws.Range("P7").Value = 0.3183113588032022
$ws.Range("S7").Value = 0.009372356013632286
$ws.Range("T7").Value = 0.009372356013632286
$ws.Range("I8").Value = 0.02944398858046029
$ws.Range("J8").Value = 0.0294439885804603
$ws.Range("M8").Value = 37.39234266666667
$ws.Range("N8").Value = 112.177028
$ws.Range("O8").Value = 0.1708807549004341
$ws.Range("P8").Value = 0.170880754900434
$ws.Range("Q8").Value = 14.15635454605911
$ws.Range("R8").Value = 127.407190914532
$ws.Range("S8").Value = 0.005031410995908815
$ws.Range("T8").Value = 0.005031410995908814
$ws.Range("I9").Value = 0.02944398858046029
$ws.Range("J9").Value = 0.0294439885804603
$ws.Range("M9").Value = 68.04127866666666
$ws.Range("N9").Value = 204.123836
$ws.Range("O9").Value = 0.3109445472994024
$ws.Range("P9").Value = 0.3109445472994024
$ws.Range("Q9").Value = 25.75972500998711
$ws.Range("R9").Value = 231.837525089884
$ws.Range("S9").Value = 0.009155447699839999
$ws.Range("T9").Value = 0.009155447699840001
$ws.Range("G10").Value = 3.441487333333333
$ws.Range("H10").Value = 10.324462
$ws.Range("I10").Value = 0.2676541983690312
$ws.Range("J10").Value = 0.2676541983690313
$ws.Range("M10").Value = 43.73434833333334
$ws.Range("N10").Value = 131.203045
$ws.Range("O10").Value = 0.1998633389969613
$ws.Range("P10").Value = 0.1998633389969613
$ws.Range("Q10").Value = 150.5112058207545
$ws.Range("R10").Value = 1354.60085238679
$ws.Range("S10").Value = 0.05349426178258961
$ws.Range("T10").Value = 0.05349426178258961
$ws.Range("G11").Value = 3.441487333333333
$ws.Range("H11").Value = 10.324462
$ws.Range("I11").Value = 0.2676541983690312
$ws.Range("J11").Value = 0.2676541983690313
$ws.Range("O11").Value = 0.3183113588032023
$ws.Range("P11").Value = 0.3183113588032022
$ws.Range("Q11").Value = 239.7109278787802
$ws.Range("R11").Value = 2157.398350909022
$ws.Range("S11").Value = 0.08519737157222816
$ws.Range("T11").Value = 0.08519737157222818
$ws.Range("G12").Value = 3.441487333333333
$ws.Range("H12").Value = 10.324462
$ws.Range("I12").Value = 0.2676541983690312
$ws.Range("J12").Value = 0.2676541983690313
$ws.Range("M12").Value = 37.39234266666667
$ws.Range("N12").Value = 112.177028
$ws.Range("O12").Value = 0.1708807549004341
$ws.Range("P12").Value = 0.170880754900434
$ws.Range("Q12").Value = 128.6852736509929
$ws.Range("R12").Value = 1158.167462858936
$ws.Range("S12").Value = 0.04573695146957057
$ws.Range("T12").Value = 0.04573695146957058
$ws.Range("G13").Value = 3.441487333333333
$ws.Range("H13").Value = 10.324462
$ws.Range("I13").Value = 0.2676541983690312
$ws.Range("J13").Value = 0.2676541983690313
$ws.Range("M13").Value = 68.04127866666666
$ws.Range("N13").Value = 204.123836
$ws.Range("O13").Value = 0.3109445472994024
$ws.Range("P13").Value = 0.3109445472994024
$ws.Range("Q13").Value = 234.1631986751368
$ws.Range("R13").Value = 2107.468788076232
$ws.Range("S13").Value = 0.08322561354464286
$ws.Range("T13").Value = 0.08322561354464288
$ws.Range("G14").Value = 2.434707333333333
$ws.Range("H14").Value = 7.304122
$ws.Range("I14").Value = 0.1893540717859783
$ws.Range("J14").Value = 0.1893540717859783
$ws.Range("M14").Value = 43.73434833333334
$ws.Range("N14").Value = 131.203045
$ws.Range("O14").Value = 0.1998633389969613
$ws.Range("P14").Value = 0.1998633389969613
$ws.Range("Q14").Value = 106.4803386057211
$ws.Range("R14").Value = 958.3230474514901
$ws.Range("S14").Value = 0.03784493703981593
$ws.Range("T14").Value = 0.03784493703981592
$ws.Range("G15").Value = 2.434707333333333
$ws.Range("H15").Value = 7.304122
$ws.Range("I15").Value = 0.1893540717859783
$ws.Range("J15").Value = 0.1893540717859783
$ws.Range("O15").Value = 0.3183113588032023
$ws.Range("P15").Value = 0.3183113588032022
$ws.Range("Q15").Value = 169.5853848810536
$ws.Range("R15").Value = 1526.268463929482
$ws.Range("S15").Value = 0.06027355188511386
$ws.Range("T15").Value = 0.06027355188511385
$ws.Range("G16").Value = 2.434707333333333
$ws.Range("H16").Value = 7.304122
$ws.Range("I16").Value = 0.1893540717859783
$ws.Range("J16").Value = 0.1893540717859783
$ws.Range("M16").Value = 37.39234266666667
$ws.Range("N16").Value = 112.177028
$ws.Range("O16").Value = 0.1708807549004341
$ws.Range("P16").Value = 0.170880754900434
$ws.Range("Q16").Value = 91.03941090104624
$ws.Range("R16").Value = 819.3546981094161
$ws.Range("S16").Value = 0.03235696673025895
$ws.Range("T16").Value = 0.03235696673025895
$ws.Range("G17").Value = 2.434707333333333
$ws.Range("H17").Value = 7.304122
$ws.Range("I17").Value = 0.1893540717859783
$ws.Range("J17").Value = 0.1893540717859783
$ws.Range("M17").Value = 68.04127866666666
$ws.Range("N17").Value = 204.123836
$ws.Range("O17").Value = 0.3109445472994024
$ws.Range("P17").Value = 0.3109445472994024
$ws.Range("Q17").Value = 165.6606001391102
$ws.Range("R17").Value = 1490.945401251992
$ws.Range("S17").Value = 0.05887861613078957
$ws.Range("T17").Value = 0.05887861613078958
